$wb = $excel.ActiveWorkbook

# The two sheets that contain these rows: "展览" (1st sheet) and "全部类型" (4th sheet)
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F4").Value = 330
    $ws.Range("F8").Value = 8344
    $ws.Range("F14").Value = 3335
    $ws.Range("F18").Value = 835
    $ws.Range("F23").Value = 1857
}
